# Insert a new data row above current row 88 (shifts rows 88-167 down to 89-168)
# and populate it with a new weekly price record for "Apio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(88).Insert()

# Columns that stay identical to the record that used to occupy row 88
# (same market / region / product / variety / quality / unit / origin / kg-or-units / classification).
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value = "Ñuble"
$ws.Cells.Item(88, 4).Value = 44574
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = 100112017
$ws.Cells.Item(88, 7).Value = "Apio"
$ws.Cells.Item(88, 8).Value = "Americana (o)"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 120
$ws.Cells.Item(88, 11).Value = 8000
$ws.Cells.Item(88, 12).Value = 8500
$ws.Cells.Item(88, 13).Value = 8250
$ws.Cells.Item(88, 14).Value = "`$/docena de matas"
$ws.Cells.Item(88, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(88, 16).Value = 1375
$ws.Cells.Item(88, 17).Value = 6
$ws.Cells.Item(88, 18).Value = "Hortaliza"

# Match the date-cell number format used by the rest of the "Fecha" column.
$ws.Cells.Item(88, 4).NumberFormat = $ws.Cells.Item(89, 4).NumberFormat
